# Insert a new record row at row 6 (pushing the existing rows 6..126 down
# to 7..127, and the sheet dimension from A1:R126 to A1:R127), then fill
# the new row with the inserted Cilantro price-record values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 6..126 down to 7..127, inheriting row 6's formatting
# (this is what gives the new D6 cell the same date style as the rest
# of column D).
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Macroferia Regional de Talca"
$ws.Range("C6").Value = "Maule"
$ws.Range("D6").Value2 = 45190
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 100112040
$ws.Range("G6").Value = "Cilantro"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 8000
$ws.Range("N6").Value = "$/caja 36 atados"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 222
$ws.Range("Q6").Value = 36
$ws.Range("R6").Value = "Hortaliza"
